$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 68.51069550467862
$ws.Range("C3").Value = 21.82653231970764
$ws.Range("C4").Value = 16.82898451760796
$ws.Range("C5").Value = 19.27929813997132
$ws.Range("C6").Value = 20.99778258680033
$ws.Range("C7").Value = 19.87336833044948
$ws.Range("C8").Value = 20.96114401019804
$ws.Range("C9").Value = 40.13926791400924
$ws.Range("C10").Value = 19.57925294315347
$ws.Range("C11").Value = 18.80933322616883
$ws.Range("C12").Value = 50.13508493338818
$ws.Range("C13").Value = 75.18997545363945
$ws.Range("C14").Value = 29.85440282233725
$ws.Range("C15").Value = 17.32348595983144
$ws.Range("C16").Value = 27.01260237527866
$ws.Range("C17").Value = 17.95003522569893
$ws.Range("C18").Value = 21.97876526869545
$ws.Range("C19").Value = 16.01351002719738
$ws.Range("C20").Value = 31.53190461125791
$ws.Range("C21").Value = 23.6246319181533
$ws.Range("C22").Value = 48.35025437515775
$ws.Range("C23").Value = 18.42986517649373
$ws.Range("C24").Value = 31.15820353891376
$ws.Range("C25").Value = 27.11105586389386
$ws.Range("C26").Value = 24.22353528893922
$ws.Range("C27").Value = 14.8590160221776
$ws.Range("C28").Value = 23.36979170070383
$ws.Range("C29").Value = 15.5782255349262
$ws.Range("C30").Value = 36.97846604231074
$ws.Range("C31").Value = 30.31328650246572
$ws.Range("C32").Value = 18.01285210148267
$ws.Range("C33").Value = 63.34770635079018
$ws.Range("C34").Value = 21.95497997601598
$ws.Range("C35").Value = 29.9528313633764
$ws.Range("C36").Value = 16.64932322838625
$ws.Range("C37").Value = 34.19079165018652
$ws.Range("C38").Value = 17.4926721642098
$ws.Range("C39").Value = 20.26149678069221
$ws.Range("C40").Value = 23.19170883719455
$ws.Range("C41").Value = 19.56377635882729
$ws.Range("C42").Value = 43.05058324576986